# Applies the numeric cell updates described by the commit diff,
# sheet by sheet, row by row (columns H-N: price/profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 15339.728
$ws.Range("I12").Value = 173.53572
$ws.Range("J12").Value = 100270.4
$ws.Range("K12").Value = 173.53572
$ws.Range("L12").Value = 100270.4
$ws.Range("M12").Value = -3.535719999999998
$ws.Range("N12").Value = -100610.4
$ws.Range("H64").Value = 48517.91
$ws.Range("I64").Value = 61629.06
$ws.Range("K64").Value = 61629.06
$ws.Range("M64").Value = -61381.06
$ws.Range("H67").Value = 48517.91
$ws.Range("I67").Value = 61629.06
$ws.Range("K67").Value = 61629.06
$ws.Range("M67").Value = -60771.06
$ws.Range("H76").Value = 3139.4167
$ws.Range("I76").Value = 3100.4285
$ws.Range("K76").Value = 3100.4285
$ws.Range("M76").Value = -2785.4285
$ws.Range("H79").Value = 3139.4167
$ws.Range("I79").Value = 3100.4285
$ws.Range("K79").Value = 3100.4285
$ws.Range("M79").Value = -2008.4285
$ws.Range("H98").Value = 1298.091
$ws.Range("I98").Value = 1484.875
$ws.Range("J98").Value = 800
$ws.Range("K98").Value = 1484.875
$ws.Range("L98").Value = 800
$ws.Range("M98").Value = 13.125
$ws.Range("N98").Value = -3796
$ws.Range("H100").Value = 2028.1333
$ws.Range("I100").Value = 2045
$ws.Range("J100").Value = 2013.375
$ws.Range("K100").Value = 2045
$ws.Range("L100").Value = 2013.375
$ws.Range("M100").Value = -1504
$ws.Range("N100").Value = -3095.375
$ws.Range("H111").Value = 10014308
$ws.Range("I111").Value = 41352.668
$ws.Range("J111").Value = 14288431
$ws.Range("K111").Value = 124058.004
$ws.Range("L111").Value = 42865293
$ws.Range("M111").Value = -120991.004
$ws.Range("N111").Value = -42871427
$ws.Range("H113").Value = 113023.336
$ws.Range("I113").Value = 144801.42
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 144801.42
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = -141547.42
$ws.Range("N113").Value = -8308
$ws.Range("H122").Value = 1298.091
$ws.Range("I122").Value = 1484.875
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 4454.625
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -2004.625
$ws.Range("N122").Value = -7300
$ws.Range("H132").Value = 16676448
$ws.Range("I132").Value = 19241440
$ws.Range("J132").Value = 4003
$ws.Range("K132").Value = 57724320
$ws.Range("L132").Value = 12009
$ws.Range("M132").Value = -57721790
$ws.Range("N132").Value = -17069
$ws.Range("H138").Value = 3494.6292
$ws.Range("I138").Value = 1071.1951
$ws.Range("J138").Value = 8226.096
$ws.Range("K138").Value = 3213.5853
$ws.Range("L138").Value = 24678.288
$ws.Range("M138").Value = 1926.4147
$ws.Range("N138").Value = -34958.288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 73683.14
$ws.Range("I2").Value = 2011.125
$ws.Range("J2").Value = 169245.83
$ws.Range("K2").Value = 2011.125
$ws.Range("L2").Value = 169245.83
$ws.Range("M2").Value = -1898.125
$ws.Range("N2").Value = -169471.83
$ws.Range("H32").Value = 35676.5
$ws.Range("I32").Value = 7215.0312
$ws.Range("J32").Value = 81214.85000000001
$ws.Range("K32").Value = 7215.0312
$ws.Range("L32").Value = 81214.85000000001
$ws.Range("M32").Value = -6928.0312
$ws.Range("N32").Value = -81788.85000000001
$ws.Range("H45").Value = 3445.375
$ws.Range("I45").Value = 2909.0667
$ws.Range("J45").Value = 4339.222
$ws.Range("K45").Value = 2909.0667
$ws.Range("L45").Value = 4339.222
$ws.Range("M45").Value = -2532.0667
$ws.Range("N45").Value = -5093.222
$ws.Range("H116").Value = 73683.14
$ws.Range("I116").Value = 2011.125
$ws.Range("J116").Value = 169245.83
$ws.Range("K116").Value = 2011.125
$ws.Range("L116").Value = 169245.83
$ws.Range("M116").Value = 282.875
$ws.Range("N116").Value = -173833.83

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 73683.14
$ws.Range("I3").Value = 2011.125
$ws.Range("J3").Value = 169245.83
$ws.Range("K3").Value = 2011.125
$ws.Range("L3").Value = 169245.83
$ws.Range("M3").Value = -1897.125
$ws.Range("N3").Value = -169473.83
$ws.Range("H107").Value = 33336882
$ws.Range("I107").Value = 83335910
$ws.Range("J107").Value = 4194.1665
$ws.Range("K107").Value = 83335910
$ws.Range("L107").Value = 4194.1665
$ws.Range("M107").Value = -83333990
$ws.Range("N107").Value = -8034.1665
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1162.75
$ws.Range("I16").Value = 750.6667
$ws.Range("J16").Value = 1410
$ws.Range("K16").Value = 750.6667
$ws.Range("L16").Value = 1410
$ws.Range("M16").Value = -463.6667
$ws.Range("N16").Value = -1984
$ws.Range("H31").Value = 1123.2
$ws.Range("I31").Value = 768.0952
$ws.Range("J31").Value = 2987.5
$ws.Range("K31").Value = 768.0952
$ws.Range("L31").Value = 2987.5
$ws.Range("M31").Value = -473.0952
$ws.Range("N31").Value = -3577.5
$ws.Range("H34").Value = 1123.2
$ws.Range("I34").Value = 768.0952
$ws.Range("J34").Value = 2987.5
$ws.Range("K34").Value = 768.0952
$ws.Range("L34").Value = 2987.5
$ws.Range("M34").Value = -566.0952
$ws.Range("N34").Value = -3391.5
$ws.Range("H113").Value = 1162.75
$ws.Range("I113").Value = 750.6667
$ws.Range("J113").Value = 1410
$ws.Range("K113").Value = 750.6667
$ws.Range("L113").Value = 1410
$ws.Range("M113").Value = 1419.3333
$ws.Range("N113").Value = -5750
$ws.Range("H132").Value = 4302.2
$ws.Range("I132").Value = 4850.3335
$ws.Range("J132").Value = 3936.7778
$ws.Range("K132").Value = 14551.0005
$ws.Range("L132").Value = 11810.3334
$ws.Range("M132").Value = -12021.0005
$ws.Range("N132").Value = -16870.3334
$ws.Range("H134").Value = 1074.8667
$ws.Range("I134").Value = 1080.2142
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 3240.6426
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -705.6425999999997
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7235.0625
$ws.Range("I5").Value = 1048.7826
$ws.Range("K5").Value = 3146.3478
$ws.Range("M5").Value = -3034.3478
$ws.Range("H37").Value = 578949.4
$ws.Range("J37").Value = 578949.4
$ws.Range("L37").Value = 1736848.2
$ws.Range("N37").Value = -1737072.2
$ws.Range("H107").Value = 818610.3
$ws.Range("J107").Value = 1070075
$ws.Range("L107").Value = 3210225
$ws.Range("N107").Value = -3214065
$ws.Range("H131").Value = 816.6022
$ws.Range("J131").Value = 890.5679
$ws.Range("L131").Value = 2671.7037
$ws.Range("N131").Value = -12751.7037
$ws.Range("H135").Value = 7235.0625
$ws.Range("I135").Value = 1048.7826
$ws.Range("K135").Value = 9439.0434
$ws.Range("M135").Value = -6904.0434
$ws.Range("H140").Value = 1729
$ws.Range("I140").Value = 1302.7858
$ws.Range("J140").Value = 2474.875
$ws.Range("K140").Value = 3908.3574
$ws.Range("L140").Value = 7424.625
$ws.Range("M140").Value = 1271.6426
$ws.Range("N140").Value = -17784.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1723.875
$ws.Range("J113").Value = 1798.5
$ws.Range("L113").Value = 1798.5
$ws.Range("N113").Value = -6138.5
$ws.Range("H132").Value = 2042.7142
$ws.Range("I132").Value = 1770.4
$ws.Range("J132").Value = 2723.5
$ws.Range("K132").Value = 5311.200000000001
$ws.Range("L132").Value = 8170.5
$ws.Range("M132").Value = -2781.200000000001
$ws.Range("N132").Value = -13230.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1648
$ws.Range("I93").Value = 1725.1428
$ws.Range("J93").Value = 1570.8572
$ws.Range("K93").Value = 1725.1428
$ws.Range("L93").Value = 1570.8572
$ws.Range("M93").Value = -477.1428000000001
$ws.Range("N93").Value = -4066.8572
$ws.Range("H116").Value = 45340.5
$ws.Range("J116").Value = 45340.5
$ws.Range("L116").Value = 45340.5
$ws.Range("N116").Value = -54518.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2320.0688
$ws.Range("I132").Value = 2076.5
$ws.Range("J132").Value = 3085.5715
$ws.Range("K132").Value = 6229.5
$ws.Range("L132").Value = 9256.7145
$ws.Range("M132").Value = -3699.5
$ws.Range("N132").Value = -14316.7145
$ws.Range("H136").Value = 1343.85
$ws.Range("I136").Value = 498.6842
$ws.Range("K136").Value = 1496.0526
$ws.Range("M136").Value = 1053.9474
